$d = $word.ActiveDocument

# --- 1. Remove the "Load/SaveProgress" and "2nd Milestone" paragraphs that
#        currently sit between "Cage needs colliders" and "3rd Milestone"
#        in the TODO section. ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Load/SaveProgress") {
        $pNext = $p.Next()
        if ($pNext.Range.Text -match "^2nd Milestone") {
            $killRange = $d.Range($p.Range.Start, $pNext.Range.End)
            $killRange.Delete()
            break
        }
    }
}

# --- 2. Add a new TODO item after "5th Milestone". ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^5th Milestone") {
        $p.Range.InsertParagraphAfter()
        $newP = $p.Next()
        $newP.Range.Text = "Disable player control (including partial) when eating/bathing/in tutorial"
        break
    }
}

# --- 3. Rename the milestone under DOING from "1st Milestone" to
#        "2nd Milestone". ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^DOING") {
        $milestoneP = $p.Next()
        $milestoneP.Range.Text = "2nd Milestone"
        break
    }
}

# --- 4. Append "Load/SaveProgress" and "1st Milestone" entries after the
#        DONE note about walk speed. ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Walk speed should be slower\.") {
        $p.Range.InsertParagraphAfter()
        $p2 = $p.Next()
        $p2.Range.Text = "Load/SaveProgress"
        $p2.Range.InsertParagraphAfter()
        $p3 = $p2.Next()
        $p3.Range.Text = "1st Milestone"
        break
    }
}
